# Apply the "danwadaag areas" update: new districts/communities appended
# to the EWEA communities/districts/FMS table (rows 184-198), with the
# first new row (184) highlighted in yellow to mark it as a new
# threshold/addition.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Build rows 184 & 189-198 from scratch (copy the formatting of the
#    last existing data row, 183, which already carries the full set of
#    per-column styles: FMS/Region bold-ish fonts, wrap-text community
#    column, small-font member column), then stamp the row height so it
#    matches the rest of the table.
#    Rows 185-188 already exist as blank placeholder rows in the sheet
#    with the right per-cell styles pre-applied, so they only need
#    values.
# ---------------------------------------------------------------------
$newRows = @(184,189,190,191,192,193,194,195,196,197,198)
foreach ($r in $newRows) {
    $ws.Range("A183:F183").Copy()
    $ws.Range("A" + $r + ":F" + $r).PasteSpecial(-4122)
    $ws.Rows.Item($r).RowHeight = 15.75
}
$excel.CutCopyMode = 0

# Row 184 is the new "header-like" highlighted row - mark it with a
# yellow fill across all six columns.
$ws.Range("A184:F184").Interior.Color = 65535

# ---------------------------------------------------------------------
# 2) Fill in the data, column by column: FMS, Region_name, District_Name,
#    Community_name, Member_name, Residence_type.
# ---------------------------------------------------------------------
$ws.Cells.Item(184,1).Value = "Southwest"
$ws.Cells.Item(184,2).Value = "Bay"
$ws.Cells.Item(184,3).Value = "Baidoa"
$ws.Cells.Item(184,4).Value = "baidoa__ca1"
$ws.Cells.Item(184,5).Value = "IOM"
$ws.Cells.Item(184,6).Value = "Urban"

$ws.Cells.Item(185,1).Value = "Southwest"
$ws.Cells.Item(185,2).Value = "Bay"
$ws.Cells.Item(185,3).Value = "Baidoa"
$ws.Cells.Item(185,4).Value = "baidoa__ca3"
$ws.Cells.Item(185,5).Value = "IOM"
$ws.Cells.Item(185,6).Value = "Urban"

$ws.Cells.Item(186,1).Value = "Southwest"
$ws.Cells.Item(186,2).Value = "Bay"
$ws.Cells.Item(186,3).Value = "Baidoa"
$ws.Cells.Item(186,4).Value = "baidoa__ca4"
$ws.Cells.Item(186,5).Value = "IOM"
$ws.Cells.Item(186,6).Value = "Urban"

$ws.Cells.Item(187,1).Value = "Southwest"
$ws.Cells.Item(187,2).Value = "Bay"
$ws.Cells.Item(187,3).Value = "Baidoa"
$ws.Cells.Item(187,4).Value = "baidoa__ca11"
$ws.Cells.Item(187,5).Value = "IOM"
$ws.Cells.Item(187,6).Value = "Urban"

$ws.Cells.Item(188,1).Value = "Southwest"
$ws.Cells.Item(188,2).Value = "Bakool"
$ws.Cells.Item(188,3).Value = "Hudur"
$ws.Cells.Item(188,4).Value = "Shida"
$ws.Cells.Item(188,5).Value = "IOM"
$ws.Cells.Item(188,6).Value = "Rural"

$ws.Cells.Item(189,1).Value = "Southwest"
$ws.Cells.Item(189,2).Value = "Bakool"
$ws.Cells.Item(189,3).Value = "Hudur"
$ws.Cells.Item(189,4).Value = "Wadajir"
$ws.Cells.Item(189,5).Value = "GREDO"
$ws.Cells.Item(189,6).Value = "Rural"

$ws.Cells.Item(190,1).Value = "Jubbaland"
$ws.Cells.Item(190,2).Value = "Lower_Juba "
$ws.Cells.Item(190,3).Value = "Kismayo"
$ws.Cells.Item(190,4).Value = "Fanole"
$ws.Cells.Item(190,5).Value = "IOM"
$ws.Cells.Item(190,6).Value = "Rural"

$ws.Cells.Item(191,1).Value = "Jubbaland"
$ws.Cells.Item(191,2).Value = "Lower_Juba "
$ws.Cells.Item(191,3).Value = "Kismayo"
$ws.Cells.Item(191,4).Value = "Calanley"
$ws.Cells.Item(191,5).Value = "NRC"
$ws.Cells.Item(191,6).Value = "Rural"

$ws.Cells.Item(192,1).Value = "Jubbaland"
$ws.Cells.Item(192,2).Value = "Lower_Juba "
$ws.Cells.Item(192,3).Value = "Kismayo"
$ws.Cells.Item(192,4).Value = "Luglow"
$ws.Cells.Item(192,5).Value = "NRC"
$ws.Cells.Item(192,6).Value = "Rural"

$ws.Cells.Item(193,1).Value = "Jubbaland"
$ws.Cells.Item(193,2).Value = "Gedo"
$ws.Cells.Item(193,3).Value = "Dolow"
$ws.Cells.Item(193,4).Value = "Qurdubey"
$ws.Cells.Item(193,5).Value = "IOM"
$ws.Cells.Item(193,6).Value = "Urban"

$ws.Cells.Item(194,1).Value = "Jubbaland"
$ws.Cells.Item(194,2).Value = "Gedo"
$ws.Cells.Item(194,3).Value = "Dolow"
$ws.Cells.Item(194,4).Value = "Kaxaarey"
$ws.Cells.Item(194,5).Value = "NRC"
$ws.Cells.Item(194,6).Value = "Urban"

$ws.Cells.Item(195,1).Value = "Southwest"
$ws.Cells.Item(195,2).Value = "Bay"
$ws.Cells.Item(195,3).Value = "Berdale"
$ws.Cells.Item(195,4).Value = "Oktober"
$ws.Cells.Item(195,5).Value = "GREDO"
$ws.Cells.Item(195,6).Value = "Urban"

$ws.Cells.Item(196,1).Value = "Jubbaland"
$ws.Cells.Item(196,2).Value = "Gedo"
$ws.Cells.Item(196,3).Value = "Baardhere"
$ws.Cells.Item(196,4).Value = "x_keskey"
$ws.Cells.Item(196,5).Value = "CONCERN"
$ws.Cells.Item(196,6).Value = "Urban"

$ws.Cells.Item(197,1).Value = "Banadir"
$ws.Cells.Item(197,2).Value = "Mogadishu"
$ws.Cells.Item(197,3).Value = "Daynile"
$ws.Cells.Item(197,4).Value = "ca19"
$ws.Cells.Item(197,5).Value = "NRC"
$ws.Cells.Item(197,6).Value = "Urban"

$ws.Cells.Item(198,1).Value = "Banadir"
$ws.Cells.Item(198,2).Value = "Mogadishu"
$ws.Cells.Item(198,3).Value = "Kahda"
$ws.Cells.Item(198,4).Value = "ca6"
$ws.Cells.Item(198,5).Value = "NRC"
$ws.Cells.Item(198,6).Value = "Urban"

# ---------------------------------------------------------------------
# 3) Update the view state: scroll the window down to the new rows and
#    leave the selection where the author left off editing.
# ---------------------------------------------------------------------
$ws.Range("A172").Select()
$excel.ActiveWindow.ScrollRow = 172
$ws.Range("F190").Select()
